$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 8062
$ws.Range("F6").Value = 299
$ws.Range("F7").Value = 807
$ws.Range("F8").Value = 617
$ws.Range("F9").Value = 101
$ws.Range("F12").Value = 872
$ws.Range("F13").Value = 3306
$ws.Range("F14").Value = 215
$ws.Range("F15").Value = 107
$ws.Range("F16").Value = 748
$ws.Range("F17").Value = 758
$ws.Range("F19").Value = 464
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 275
$ws.Range("F22").Value = 236
$ws.Range("F23").Value = 348
$ws.Range("F26").Value = 125
$ws.Range("F27").Value = 288
$ws.Range("F28").Value = 32
$ws.Range("F32").Value = 569
$ws.Range("F33").Value = 27
$ws.Range("F35").Value = 17
$ws.Range("F36").Value = 23
$ws.Range("F38").Value = 107

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 209

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 209
$ws.Range("F5").Value = 8062
$ws.Range("F8").Value = 299
$ws.Range("F9").Value = 807
$ws.Range("F10").Value = 617
$ws.Range("F11").Value = 101
$ws.Range("F14").Value = 872
$ws.Range("F16").Value = 3307
$ws.Range("F17").Value = 215
$ws.Range("F18").Value = 107
$ws.Range("F20").Value = 748
$ws.Range("F21").Value = 758
$ws.Range("F24").Value = 464
$ws.Range("F25").Value = 23
$ws.Range("F26").Value = 275
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 349
$ws.Range("F30").Value = 133
$ws.Range("F31").Value = 125
$ws.Range("F32").Value = 288
$ws.Range("F33").Value = 32
$ws.Range("F37").Value = 569
$ws.Range("F38").Value = 27
$ws.Range("F40").Value = 17
$ws.Range("F41").Value = 23
$ws.Range("F43").Value = 107

